$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-14 Friday" "2025-03-15 Saturday"

Replace-Text "22×73=1606" "69×46=3174"
Replace-Text "14×19=266" "25×64=1600"
Replace-Text "13×49=637" "83×14=1162"
Replace-Text "93×48=4464" "30×88=2640"
Replace-Text "23×22=506" "94×64=6016"

Replace-Text "59×92=5428" "63×96=6048"
Replace-Text "53×88=4664" "23×45=1035"
Replace-Text "13×79=1027" "72×76=5472"
Replace-Text "63×98=6174" "87×60=5220"
Replace-Text "13×39=507" "82×67=5494"

Replace-Text "91×29=2639" "62×14=868"
Replace-Text "12×91=1092" "81×52=4212"
Replace-Text "15×35=525" "89×23=2047"
Replace-Text "60×69=4140" "71×28=1988"
Replace-Text "25×22=550" "20×52=1040"

Replace-Text "74×11=814" "19×73=1387"
Replace-Text "88×30=2640" "91×94=8554"
Replace-Text "25×17=425" "11×66=726"
Replace-Text "49×67=3283" "57×45=2565"
Replace-Text "34×31=1054" "72×18=1296"

Replace-Text "88×69=6072" "65×62=4030"
Replace-Text "81×78=6318" "65×31=2015"
Replace-Text "33×45=1485" "15×29=435"
Replace-Text "89×90=8010" "27×53=1431"
Replace-Text "16×48=768" "84×15=1260"
